$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (B1:I1) with the new class-boundary labels,
# replacing the old "min, SD_nedre, SD_D, D_M, M_G, G_SG, SG_øvre, max"
$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"

# Move/update the active selection to A2, as reflected in the sheetView
$ws.Range("A2").Select()
